# The workbook originally contains a single sheet named "SheetA".
# This edit:
#   1. Duplicates "SheetA" (same data/layout) right after itself.
#   2. Renames the original sheet to "Monkey" (stays first, keeps sheetId 1).
#   3. Names the duplicate "SheetA" (becomes second, new sheetId 2).
#   4. Makes the new "SheetA" the active/selected tab (activeTab=1),
#      while "Monkey" is no longer the selected tab.
#   5. Sets the selection on both sheets to the single cell A20
#      (previously the range A20:A25).

$wb = $excel.ActiveWorkbook

$original = $wb.Worksheets.Item("SheetA")

# Duplicate the sheet, placing the copy right after the original.
$original.Copy([System.Reflection.Missing]::Value, $original)

# The copy is now immediately after $original in the tab order.
$duplicate = $wb.Worksheets.Item($original.Index + 1)

# Rename: original -> "Monkey", duplicate -> "SheetA"
$duplicate.Name = "SheetA_tmp_rename"
$original.Name = "Monkey"
$duplicate.Name = "SheetA"

# Fix up the selection on each sheet to a single cell (A20) and make the
# new "SheetA" the active sheet/tab.
$original.Activate()
$original.Range("A20").Select() | Out-Null
$duplicate.Activate()
$duplicate.Range("A20").Select() | Out-Null
